$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update priority values in column A for the affected backlog rows
$ws.Range("A7").Value = 0
$ws.Range("A8").Value = 0
$ws.Range("A12").Value = 1
$ws.Range("A15").Value = 1

# Update the active selection on the sheet
$ws.Range("A11").Select()
